# Auto-generated: updates currentAveragePrice / LevePrice / LeveProfit columns
# across several Leve tables, refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 71431930
$ws.Range("I62").Value = 166669460
$ws.Range("J62").Value = 3788.25
$ws.Range("K62").Value = 166669460
$ws.Range("L62").Value = 3788.25
$ws.Range("M62").Value = -166668836
$ws.Range("N62").Value = -5036.25
$ws.Range("H65").Value = 71431930
$ws.Range("I65").Value = 166669460
$ws.Range("J65").Value = 3788.25
$ws.Range("K65").Value = 833347300
$ws.Range("L65").Value = 18941.25
$ws.Range("M65").Value = -833344180
$ws.Range("N65").Value = -25181.25
$ws.Range("H86").Value = 2810.3215
$ws.Range("I86").Value = 2341.182
$ws.Range("J86").Value = 3113.8823
$ws.Range("K86").Value = 2341.182
$ws.Range("L86").Value = 3113.8823
$ws.Range("M86").Value = -1218.182
$ws.Range("N86").Value = -5359.8823
$ws.Range("H89").Value = 2810.3215
$ws.Range("I89").Value = 2341.182
$ws.Range("J89").Value = 3113.8823
$ws.Range("K89").Value = 11705.91
$ws.Range("L89").Value = 15569.4115
$ws.Range("M89").Value = -6089.91
$ws.Range("N89").Value = -26801.4115
$ws.Range("H129").Value = 876.44446
$ws.Range("I129").Value = 517.75
$ws.Range("J129").Value = 1163.4
$ws.Range("K129").Value = 1553.25
$ws.Range("L129").Value = 3490.2
$ws.Range("M129").Value = 3446.75
$ws.Range("N129").Value = -13490.2
$ws.Range("H131").Value = 47621340
$ws.Range("I131").Value = 90910750
$ws.Range("J131").Value = 2990
$ws.Range("K131").Value = 272732250
$ws.Range("L131").Value = 8970
$ws.Range("M131").Value = -272727210
$ws.Range("N131").Value = -19050
$ws.Range("H138").Value = 13295.966
$ws.Range("I138").Value = 1161.8148
$ws.Range("J138").Value = 33151.848
$ws.Range("K138").Value = 3485.4444
$ws.Range("L138").Value = 99455.54399999999
$ws.Range("M138").Value = 1654.5556
$ws.Range("N138").Value = -109735.544

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 27448.334
$ws.Range("J24").Value = 27448.334
$ws.Range("L24").Value = 27448.334
$ws.Range("N24").Value = -28196.334
$ws.Range("H61").Value = 1823.9231
$ws.Range("I61").Value = 1634.5555
$ws.Range("K61").Value = 1634.5555
$ws.Range("M61").Value = -1422.5555
$ws.Range("H100").Value = 27448.334
$ws.Range("J100").Value = 27448.334
$ws.Range("L100").Value = 27448.334
$ws.Range("N100").Value = -29612.334
$ws.Range("H136").Value = 1823.9231
$ws.Range("I136").Value = 1634.5555
$ws.Range("K136").Value = 4903.666499999999
$ws.Range("M136").Value = -2353.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1346.16
$ws.Range("I94").Value = 1003
$ws.Range("J94").Value = 1860.9
$ws.Range("K94").Value = 1003
$ws.Range("L94").Value = 1860.9
$ws.Range("M94").Value = -552
$ws.Range("N94").Value = -2762.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2833.8
$ws.Range("I134").Value = 1731.75
$ws.Range("J134").Value = 4486.875
$ws.Range("K134").Value = 5195.25
$ws.Range("L134").Value = 13460.625
$ws.Range("M134").Value = -2660.25
$ws.Range("N134").Value = -18530.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6536740
$ws.Range("I102").Value = 7937290.5
$ws.Range("J102").Value = 838
$ws.Range("K102").Value = 7937290.5
$ws.Range("L102").Value = 838
$ws.Range("M102").Value = -7935668.5
$ws.Range("N102").Value = -4082
$ws.Range("H122").Value = 1084.5
$ws.Range("I122").Value = 951.75
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 2855.25
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -405.25
$ws.Range("N122").Value = -8950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1514.0834
$ws.Range("J7").Value = 1775
$ws.Range("L7").Value = 1775
$ws.Range("N7").Value = -1999
$ws.Range("H40").Value = 2851.647
$ws.Range("I40").Value = 2885.2
$ws.Range("J40").Value = 2600
$ws.Range("K40").Value = 2885.2
$ws.Range("L40").Value = 2600
$ws.Range("M40").Value = -2749.2
$ws.Range("N40").Value = -2872
$ws.Range("H55").Value = 473.54285
$ws.Range("I55").Value = 250.6923
$ws.Range("J55").Value = 605.2273
$ws.Range("K55").Value = 250.6923
$ws.Range("L55").Value = 605.2273
$ws.Range("M55").Value = -77.69229999999999
$ws.Range("N55").Value = -951.2273
$ws.Range("H61").Value = 1520.3334
$ws.Range("I61").Value = 1152
$ws.Range("J61").Value = 1888.6666
$ws.Range("K61").Value = 1152
$ws.Range("L61").Value = 1888.6666
$ws.Range("M61").Value = -950
$ws.Range("N61").Value = -2292.6666
$ws.Range("H113").Value = 1520.3334
$ws.Range("I113").Value = 1152
$ws.Range("J113").Value = 1888.6666
$ws.Range("K113").Value = 1152
$ws.Range("L113").Value = 1888.6666
$ws.Range("M113").Value = 1018
$ws.Range("N113").Value = -6228.6666
$ws.Range("H122").Value = 2257.3157
$ws.Range("I122").Value = 1754.3334
$ws.Range("J122").Value = 2710
$ws.Range("K122").Value = 5263.0002
$ws.Range("L122").Value = 8130
$ws.Range("M122").Value = -2813.0002
$ws.Range("N122").Value = -13030
$ws.Range("H126").Value = 1514.0834
$ws.Range("J126").Value = 1775
$ws.Range("L126").Value = 5325
$ws.Range("N126").Value = -10265
$ws.Range("H136").Value = 5960.8184
$ws.Range("I136").Value = 1494.5
$ws.Range("J136").Value = 13776.875
$ws.Range("K136").Value = 4483.5
$ws.Range("L136").Value = 41330.625
$ws.Range("M136").Value = -1933.5
$ws.Range("N136").Value = -46430.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 31983.334
$ws.Range("I2").Value = 60966.668
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 60966.668
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -60854.668
$ws.Range("N2").Value = -3224
$ws.Range("H81").Value = 1221.2858
$ws.Range("I81").Value = 1275
$ws.Range("J81").Value = 1199.8
$ws.Range("K81").Value = 2550
$ws.Range("L81").Value = 2399.6
$ws.Range("M81").Value = -1489
$ws.Range("N81").Value = -4521.6
$ws.Range("H84").Value = 1221.2858
$ws.Range("I84").Value = 1275
$ws.Range("J84").Value = 1199.8
$ws.Range("K84").Value = 12750
$ws.Range("L84").Value = 11998
$ws.Range("M84").Value = -7446
$ws.Range("N84").Value = -22606
$ws.Range("H126").Value = 1236
$ws.Range("I126").Value = 1236
$ws.Range("K126").Value = 3708
$ws.Range("M126").Value = -1238
$ws.Range("H136").Value = 5103.7803
$ws.Range("I136").Value = 2841.6
$ws.Range("J136").Value = 8638.4375
$ws.Range("K136").Value = 8524.799999999999
$ws.Range("L136").Value = 25915.3125
$ws.Range("M136").Value = -5974.799999999999
$ws.Range("N136").Value = -31015.3125
